$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name to reflect new "through" date
$ws.Name = "Through 2022-04-25"

# Update header label in I1 (shared string) to match new "through" date
$ws.Range("I1").Value = "2022 (through 04-25)"

# Update data values for May (row 5) and Total (row 14) in column I
$ws.Range("I5").Value = 106
$ws.Range("I14").Value = 541
